$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifting old B -> C)
$ws.Range("B1").EntireColumn.Insert()

# Row 1 header
$ws.Range("B1").Value = "Trening"

# Row 2 (existing 10-15 row): add Trening value, update Acceleration_SMA value
$ws.Range("B2").Value = "Duża Gra"
$ws.Range("C2").Value = 3.471837611476191

# Insert new row 3 for 10-15 / Mała Gra
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A3").Value = "10-15"
$ws.Range("B3").Value = "Mała Gra"
$ws.Range("C3").Value = 2.408517088238095

# Row 4 (originally row 3, the "5-10" row) now needs Trening + updated Acceleration_SMA
$ws.Range("B4").Value = "Duża Gra"
$ws.Range("C4").Value = 3.200018394523809

# Insert new row 5 for 5-10 / Mała Gra
$ws.Range("A5").EntireRow.Insert()
$ws.Range("A5").Value = "5-10"
$ws.Range("B5").Value = "Mała Gra"
$ws.Range("C5").Value = 2.51308704
